$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.601.59'
$ws.Range("E2").Value = '  +0.33%  '

$ws.Range("D3").Value = '3.149.79'
$ws.Range("E3").Value = '  +0.08%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '''528.97'
$ws.Range("E5").Value = '  -1.36%  '

$ws.Range("D6").Value = '''139.28'
$ws.Range("E6").Value = '  -0.43%  '

$ws.Range("D8").Value = '''0.535'
$ws.Range("E8").Value = '  +14.00%  '

$ws.Range("D9").Value = '''7.30'
$ws.Range("E9").Value = '  -0.20%  '

$ws.Range("E10").Value = '  +5.34%  '

$ws.Range("D11").Value = '''0.111'
$ws.Range("E11").Value = '  +3.05%  '

$ws.Range("E12").Value = '  +3.20%  '

$ws.Range("D13").Value = '3.699.77'
$ws.Range("E13").Value = '  +1.00%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '''25.75'
$ws.Range("E14").Value = '  +0.04%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '''0.0000172'
$ws.Range("E15").Value = '  +3.98%  '

$ws.Range("D16").Value = '58.681.66'
$ws.Range("E16").Value = '  +0.36%  '

$ws.Range("D17").Value = '''6.25'
$ws.Range("E17").Value = '  +2.84%  '

$ws.Range("D18").Value = '3.160.40'
$ws.Range("E18").Value = '  +0.87%  '

$ws.Range("D19").Value = '''12.94'
$ws.Range("E19").Value = '  +1.49%  '

$ws.Range("D20").Value = '''8.12'
$ws.Range("E20").Value = '  -0.88%  '

$ws.Range("D21").Value = '''372.50'
$ws.Range("E21").Value = '  +3.13%  '

$ws.Range("E22").Value = '  +1.69%  '

$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("E24").Value = '  +4.17%  '

$ws.Range("D25").Value = '''69.46'
$ws.Range("E25").Value = '  +0.26%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").Value = '''0.999'
$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("D28").Value = '''8.28'
$ws.Range("E28").Value = '  +12.43%  '

$ws.Range("D29").Value = '0.0₃0862'
$ws.Range("E29").Value = '  -2.51%  '

$ws.Range("D30").Value = '''22.22'
$ws.Range("E30").Value = '  +2.27%  '

$ws.Range("E31").Value = '  -0.56%  '

$ws.Range("D32").Value = '''6.06'
$ws.Range("E32").Value = '  -1.91%  '

$ws.Range("E33").Value = '  -0.62%  '

$ws.Range("D34").Value = '''1.15'
$ws.Range("E34").Value = '  -0.41%  '

$ws.Range("D35").Value = '''6.28'
$ws.Range("E35").Value = '  +2.56%  '

$ws.Range("D36").Value = '''158.13'
$ws.Range("E36").Value = '  -0.71%  '

$ws.Range("D37").Value = '''1.33'
$ws.Range("E37").Value = '  +4.64%  '

$ws.Range("D38").Value = '''24.99'
$ws.Range("E38").Value = '  -3.65%  '

$ws.Range("D39").Value = '''1.67'
$ws.Range("E39").Value = '  -1.53%  '

$ws.Range("D40").Value = '''0.0684'
$ws.Range("E40").Value = '  +1.20%  '

$ws.Range("D41").Value = '2.625.85'
$ws.Range("E41").Value = '  +4.88%  '

$ws.Range("E42").Value = '  +4.76%  '

$ws.Range("D43").Value = '''0.718'
$ws.Range("E43").Value = '  +1.90%  '

$ws.Range("D44").Value = '''39.01'
$ws.Range("E44").Value = '  +3.73%  '

$ws.Range("E45").Value = '  +5.99%  '

$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("D47").Value = '3.196.29'
$ws.Range("E47").Value = '  +0.38%  '

$ws.Range("E48").Value = '  +12.78%  '

$ws.Range("E49").Value = '  +1.54%  '

$ws.Range("D50").Value = '''0.976'
$ws.Range("E50").Value = '  -2.44%  '

$ws.Range("D51").Value = '''20.03'
$ws.Range("E51").Value = '  +0.14%  '
